$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Supuestos" (column J) values that were missing on the
# originally existing rows 2-7.
$ws.Range("J2").Value = "aaaa"
$ws.Range("J3").Value = "aasdfdfgh"
$ws.Range("J4").Value = "aasdfdfgh"
$ws.Range("J5").Value = "aaa"
$ws.Range("J6").Value = "Este es un supuesto"
$ws.Range("J7").Value = "aa"

# Insert a brand-new record just above the current row 8, pushing the old
# rows 8-12 down to become rows 9-13.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "Gobierno Ciudadano`n"
$ws.Range("B8").Value = "Secretaría de Finanzas"
$ws.Range("C8").Value = 2022
$ws.Range("D8").Value = "Componente"
$ws.Range("E8").Value = 5529
$ws.Range("F8").Value = "Programa Presupuestario 3"
$ws.Range("G8").Value = "Optimizar la generación de recursos propios estatales "
$ws.Range("H8").Value = "Subasta pública electrónica"
# "2" must land as text (matching the source report's inline-string cells),
# not auto-converted to a number. Force text entry, then restore the
# original (default) cell style so no extra style slot is introduced.
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "2"
$ws.Range("I8").Style = $ws.Range("I9").Style
$ws.Range("J8").Value = "Supuesto de la actividad Aprovechar los recursos propios"
$ws.Range("K8").Value = "Area 1"

# Fill in the "Supuestos" (column J) values for the rows that were shifted
# down from their original positions (old rows 8-12, now rows 9-13).
$ws.Range("J9").Value = "qqqq"
$ws.Range("J10").Value = "aaaaaa"
$ws.Range("J11").Value = "ggg"
$ws.Range("J12").Value = "TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst TEst vv"
$ws.Range("J13").Value = "aaaa"
